$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the shared "ping files" comment text from (212) to (215) everywhere it occurs.
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*shaconemo (212) ping files*") {
        $cell.Value2 = $val -replace "shaconemo \(212\) ping files", "shaconemo (215) ping files"
    }
}

# 2. Add the forgotten NEMO code name comments for rows 548-552 (expfe, expsi,
#    expcalc, ppdiat, ppmisc) and refresh the comment-author cell in column F.
$ws.Cells.Item(548, 5).Value2 = "Not available in LPJ-GUESS.  Available in PISCES: expfe"
$ws.Cells.Item(549, 5).Value2 = "Not available in LPJ-GUESS.  Available in PISCES: expsi"
$ws.Cells.Item(550, 5).Value2 = "Not available in LPJ-GUESS.  Available in PISCES: expcal"
$ws.Cells.Item(551, 5).Value2 = "Not available in LPJ-GUESS.  Available in PISCES: PPPHY2"
$ws.Cells.Item(552, 5).Value2 = "Not available in LPJ-GUESS.  Available in PISCES: PPPHY2"

for ($r = 548; $r -le 552; $r++) {
    $ws.Cells.Item($r, 6).Value2 = "David Warlind, Raffaele Bernardello"
}

# 3. Leave the view the way the author left it: scrolled back to the top with
#    the last-edited cell (E552) selected.
[void]$ws.Range("E552").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
